$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: update coordinates in F10 and set G10 validation label
$ws.Range("F10").Value = "74.02583599090578,30.05605700974394,0 74.02549266815187,30.05323396722543,0 74.02433395385744,30.049705050900876,0 74.02364730834962,30.042869739559226,0 74.04948234558107,30.03870888434231,0 74.05317306518556,30.049147842091724,0 74.05381679534914,30.05427404488543,0 74.04815196990968,30.056874191241224,0 74.04411792755128,30.058359958502557,0 74.04089927673341,30.058879971776534,0 74.03677940368654,30.060142849784942,0 74.0272092819214,30.06170285331865,0"
$ws.Range("G10").Value = "Bad not fixable"

# Row 14: set G14 validation label
$ws.Range("G14").Value = "Good Polygon"
